$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: new titrate sample data
$ws.Range("A23").NumberFormat = "m/d/yy"
$ws.Range("A23").Value = 43215

$ws.Range("B23").Value = 2218.2645308156002

$ws.Range("F23").Value = "end of sample"

# D23 already carries the shared formula; re-set it so it recalculates against
# the newly populated B23/C23
$ws.Range("D23").Formula = "=100*(B23-C23)/C23"

# Row 24: A24 gets the same date style as the rest of column A but remains blank
$ws.Range("A24").NumberFormat = "m/d/yy"

# Update the active selection to match
$ws.Range("F24").Select()
